$wb = $excel.ActiveWorkbook

# --- Remove Sheet2 ---------------------------------------------------------
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Sheet2").Delete()
$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename header row (shared strings) ------------------------------------
$ws.Range("A1").Value = "Kolona 1"
$ws.Range("B1").Value = "Kolona 2"
$ws.Range("C1").Value = "Kolona 3"
$ws.Range("D1").Value = "Kolona 4"

# --- New data row ------------------------------------------------------------
$ws.Range("A2").Value = "Vrijednost "
$ws.Range("B2").Value = 12345
$ws.Range("C2").Value = "nerminsehic1993@gmail.com"
$ws.Range("D2").Value = "bla bla 4"

# --- Hyperlink on C2 (adds Hyperlink font/style automatically) -------------
[void]$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:nerminsehic1993@gmail.com")

# --- Column C width ----------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 25.5

# --- Selection / active sheet -------------------------------------------------
[void]$ws.Select()
[void]$ws.Range("E9").Select()
